$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B values for rows 2-8 (scaled / recomputed values)
$ws.Range("B2").Value = 116.1461029052734
$ws.Range("B3").Value = 90.70700073242188
$ws.Range("B4").Value = 30.65430068969727
$ws.Range("B5").Value = -1.531499981880188
$ws.Range("B6").Value = -75.86190032958984
$ws.Range("B7").Value = -60.91680145263672
$ws.Range("B8").Value = 85.29180145263672

# Add new row 9 - "Пастбище" row
$ws.Range("A9").Value = "Пастбище"
$ws.Range("B9").Value = 184.489013671875
$ws.Range("C9").Value = 0.0241
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.4611000120639801
$ws.Range("F9").Value = -0.06889999657869339
